$d = $word.ActiveDocument

# Step 1: change " crude odds ratio" -> " crude odds ratio, " (append ", ")
$d.Content.Find.Execute(" crude odds ratio", $true, $false, $false, $false, $false, $true, 1, $false, " crude odds ratio, ", 2)

# Step 2: locate the run we just edited and collapse to its end
$rng = $d.Content
$rng.Find.Execute(" crude odds ratio, ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(0)

# Step 3: insert the new abbreviation text (inherits the non-italic formatting
# of the preceding run), then italicize only the "AIC" portion.
$insertStart = $rng.Start
$rng.InsertAfter("AIC Akaike Information Criterion")

$italicRng = $d.Range($insertStart, $insertStart + 3)
$italicRng.Font.Italic = $true
